$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM (transcripts-per-million) derived NATMI ligand-receptor
# edge table for Il12a -> Il12rb1: values recomputed, plus a new set of
# rows (6-9) for the MuSCs sending cluster that did not exist before.

$row2 = @("ECs", "Il12a", "Il12rb1", "ECs", 1, 0.3333333333333333, 0.019271, 0.057813, 0.009412517461340836, 0.009412517461340836, 3, 1, 0.569907, 1.709721, 0.1475417602688563, 0.1475417602688563, 0.010982677797, 0.09884410017300001, 0.001388739394807574, 0.001388739394807574)
for ($c = 0; $c -lt $row2.Length; $c++) {
    $ws.Cells.Item(2, $c + 1).Value = $row2[$c]
}

$row3 = @("ECs", "Il12a", "Il12rb1", "FAPs", 1, 0.3333333333333333, 0.019271, 0.057813, 0.009412517461340836, 0.009412517461340836, 3, 1, 1.694965, 5.084895, 0.4388051378454766, 0.4388051378454766, 0.032663670515, 0.293973034635, 0.004130261022096621, 0.004130261022096621)
for ($c = 0; $c -lt $row3.Length; $c++) {
    $ws.Cells.Item(3, $c + 1).Value = $row3[$c]
}

$row4 = @("ECs", "Il12a", "Il12rb1", "MuSCs", 1, 0.3333333333333333, 0.019271, 0.057813, 0.009412517461340836, 0.009412517461340836, 2, 0.6666666666666666, 0.125128, 0.375384, 0.03239406671425592, 0.03239406671425593, 0.002411341688, 0.021702075192, 0.0003049097185917738, 0.0003049097185917739)
for ($c = 0; $c -lt $row4.Length; $c++) {
    $ws.Cells.Item(4, $c + 1).Value = $row4[$c]
}

$row5 = @("ECs", "Il12a", "Il12rb1", "Resolving-Mac", 1, 0.3333333333333333, 0.019271, 0.057813, 0.009412517461340836, 0.009412517461340836, 3, 1, 1.472682666666667, 4.418048, 0.3812590351714111, 0.3812590351714111, 0.02838006766933333, 0.255420609024, 0.003588607325844866, 0.003588607325844866)
for ($c = 0; $c -lt $row5.Length; $c++) {
    $ws.Cells.Item(5, $c + 1).Value = $row5[$c]
}

$row6 = @("MuSCs", "Il12a", "Il12rb1", "ECs", 3, 1, 2.028109, 6.084327, 0.9905874825386592, 0.9905874825386591, 3, 1, 0.569907, 1.709721, 0.1475417602688563, 0.1475417602688563, 1.155833515863, 10.402501642767, 0.1461530208740488, 0.1461530208740487)
for ($c = 0; $c -lt $row6.Length; $c++) {
    $ws.Cells.Item(6, $c + 1).Value = $row6[$c]
}

$row7 = @("MuSCs", "Il12a", "Il12rb1", "FAPs", 3, 1, 2.028109, 6.084327, 0.9905874825386592, 0.9905874825386591, 3, 1, 1.694965, 5.084895, 0.4388051378454766, 0.4388051378454766, 3.437573771185, 30.938163940665, 0.4346748768233801, 0.43467487682338)
for ($c = 0; $c -lt $row7.Length; $c++) {
    $ws.Cells.Item(7, $c + 1).Value = $row7[$c]
}

$row8 = @("MuSCs", "Il12a", "Il12rb1", "MuSCs", 3, 1, 2.028109, 6.084327, 0.9905874825386592, 0.9905874825386591, 2, 0.6666666666666666, 0.125128, 0.375384, 0.03239406671425592, 0.03239406671425593, 0.253773222952, 2.283959006568, 0.03208915699566416, 0.03208915699566416)
for ($c = 0; $c -lt $row8.Length; $c++) {
    $ws.Cells.Item(8, $c + 1).Value = $row8[$c]
}

$row9 = @("MuSCs", "Il12a", "Il12rb1", "Resolving-Mac", 3, 1, 2.028109, 6.084327, 0.9905874825386592, 0.9905874825386591, 3, 1, 1.472682666666667, 4.418048, 0.3812590351714111, 0.3812590351714111, 2.986760970410667, 26.880848733696, 0.3776704278455663, 0.3776704278455662)
for ($c = 0; $c -lt $row9.Length; $c++) {
    $ws.Cells.Item(9, $c + 1).Value = $row9[$c]
}
